# Generate Report for Handoff
# Updates status/priority/handoff-datetime for the file
# "6c7dd511-eaf1-4059-9cec-38719905fb56.md" which moved from
# "In Translation" to "Ready for handoff", and refreshes the handoff
# timestamps for the files already in the handoff batch
# (b470d857-... and ba8b2921-...).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ovw = $wb.Worksheets.Item("Overview")
# Row 6 -> 6c7dd511-eaf1-4059-9cec-38719905fb56.md
$ovw.Range("E6").Value = "Ready for handoff"
$ovw.Range("F6").Value = "Ready for handoff"
$ovw.Range("G6").Value = "2016-09-07 16:29:31"
# Row 7 -> b470d857-e8d4-43c9-8572-c627035024fb.md
$ovw.Range("G7").Value = "2016-09-07 16:29:31"
# Row 8 -> ba8b2921-856d-4ce8-b2fb-39d9aea92fe3.md
$ovw.Range("G8").Value = "2016-09-07 16:29:31"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
# Row 6 -> 6c7dd511-eaf1-4059-9cec-38719905fb56.md
$zh.Range("C6").Value = "Ready for handoff"
$zh.Range("E6").Value = "ht"
$zh.Range("H6").Value = "2016-09-07 16:29:25"
# Row 7 -> b470d857-e8d4-43c9-8572-c627035024fb.md
$zh.Range("E7").Value = "ht"
$zh.Range("H7").Value = "2016-09-07 16:29:25"
# Row 8 -> ba8b2921-856d-4ce8-b2fb-39d9aea92fe3.md
$zh.Range("E8").Value = "ht"
$zh.Range("H8").Value = "2016-09-07 16:29:25"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
# Row 6 -> 6c7dd511-eaf1-4059-9cec-38719905fb56.md
$de.Range("C6").Value = "Ready for handoff"
$de.Range("E6").Value = "ht"
$de.Range("H6").Value = "2016-09-07 16:29:31"
# Row 7 -> b470d857-e8d4-43c9-8572-c627035024fb.md
$de.Range("E7").Value = "ht"
$de.Range("H7").Value = "2016-09-07 16:29:31"
# Row 8 -> ba8b2921-856d-4ce8-b2fb-39d9aea92fe3.md
$de.Range("E8").Value = "ht"
$de.Range("H8").Value = "2016-09-07 16:29:31"
